$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (text cells) ---
# A2: "fcff" -> "ff" (plain text change, safe to assign directly)
$ws.Range("A2").Value = "ff"

# D2: "0" -> "1" (must stay text; copy from K2 which already holds text "1"
# so the cell keeps its string type/no extra number formatting)
$ws.Range("K2").Copy($ws.Range("D2"))

# F2: "1" -> "0" (must stay text; copy from C2 which already holds text "0")
$ws.Range("C2").Copy($ws.Range("F2"))

# --- Row 3 (numeric summary cells) ---
$ws.Range("D3").Value = 1
$ws.Range("F3").Value = 0
